$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price strings remain text (matches source formatting)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "57.173.19"
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("D3").Value = "2.513.97"
$ws.Range("E3").Value = "  +2.26%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "495.54"
$ws.Range("E5").Value = "  +3.64%  "
$ws.Range("D6").Value = "153.40"
$ws.Range("E6").Value = "  +10.34%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "0.517"
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("D9").Value = "2.535.78"
$ws.Range("E9").Value = "  +3.19%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  +5.43%  "
$ws.Range("E11").Value = "  +5.75%  "
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  +5.09%  "
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("D14").Value = "2.955.63"
$ws.Range("E14").Value = "  +3.46%  "
$ws.Range("D15").Value = "57.333.96"
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("D16").Value = "21.36"
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("D17").Value = "0.0000138"
$ws.Range("E17").Value = "  +3.58%  "
$ws.Range("D18").Value = "2.536.18"
$ws.Range("E18").Value = "  +3.16%  "
$ws.Range("E19").Value = "  +5.92%  "
$ws.Range("D20").Value = "10.34"
$ws.Range("E20").Value = "  +4.16%  "
$ws.Range("D21").Value = "324.25"
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("D22").Value = "5.99"
$ws.Range("E22").Value = "  +6.55%  "
$ws.Range("D23").Value = "0.997"
$ws.Range("E23").Value = "  +0.35%  "
$ws.Range("D24").Value = "58.58"
$ws.Range("E24").Value = "  +2.58%  "
$ws.Range("D25").Value = "0.411"
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("E27").Value = "  -0.40%  "
$ws.Range("D28").Value = "2.622.95"
$ws.Range("E28").Value = "  +2.67%  "
$ws.Range("D29").Value = "7.65"
$ws.Range("E29").Value = "  +4.58%  "
$ws.Range("D30").Value = "0.0₃0830"
$ws.Range("E30").Value = "  +7.94%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "151.82"
$ws.Range("E32").Value = "  +2.09%  "
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("D34").Value = "18.43"
$ws.Range("E34").Value = "  +3.03%  "
$ws.Range("E35").Value = "  +3.09%  "
$ws.Range("D36").Value = "0.910"
$ws.Range("E36").Value = "  +7.56%  "
$ws.Range("D37").Value = "3.85"
$ws.Range("E37").Value = "  +7.36%  "
$ws.Range("E38").Value = "  +4.64%  "
$ws.Range("E39").Value = "  +10.90%  "
$ws.Range("D40").Value = "34.43"
$ws.Range("E40").Value = "  +3.23%  "
$ws.Range("D41").Value = "3.54"
$ws.Range("E41").Value = "  +4.48%  "
$ws.Range("D42").Value = "0.620"
$ws.Range("E42").Value = "  +3.06%  "
$ws.Range("D43").Value = "0.0564"
$ws.Range("E43").Value = "  +4.28%  "
$ws.Range("D44").Value = "0.995"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("D45").Value = "4.94"
$ws.Range("E45").Value = "  +6.07%  "
$ws.Range("D46").Value = "268.70"
$ws.Range("E46").Value = "  +5.16%  "
$ws.Range("D47").Value = "0.0949"
$ws.Range("E47").Value = "  +6.31%  "
$ws.Range("D48").Value = "0.0230"
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("D50").Value = "18.16"
$ws.Range("E50").Value = "  +6.66%  "
$ws.Range("D51").Value = "1.904.61"
$ws.Range("E51").Value = "  -1.13%  "
